# Auto-generated Excel COM-interop script
# Applies updated market-price figures (currentAveragePrice / NQ / HQ / Leve profit columns)
# to the leve-profit tracking sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 613.625
$ws.Range("I33").Value = 401.5
$ws.Range("K33").Value = 401.5
$ws.Range("M33").Value = -172.5

$ws.Range("H40").Value = 6349.2
$ws.Range("I40").Value = 6349.2
$ws.Range("K40").Value = 6349.2
$ws.Range("M40").Value = -6174.2

$ws.Range("H62").Value = 18707.438
$ws.Range("J62").Value = 17180.223
$ws.Range("L62").Value = 17180.223
$ws.Range("N62").Value = -18428.223

$ws.Range("H65").Value = 18707.438
$ws.Range("J65").Value = 17180.223
$ws.Range("L65").Value = 85901.11500000001
$ws.Range("N65").Value = -92141.11500000001

$ws.Range("H96").Value = 366.18182
$ws.Range("I96").Value = 377.77777
$ws.Range("K96").Value = 1133.33331
$ws.Range("M96").Value = 239.66669

$ws.Range("H132").Value = 1539.2778
$ws.Range("I132").Value = 1286.6863
$ws.Range("K132").Value = 3860.0589
$ws.Range("M132").Value = -1330.0589

$ws.Range("H137").Value = 1647962.2
$ws.Range("I137").Value = 1827.7142
$ws.Range("J137").Value = 2471029.5
$ws.Range("K137").Value = 5483.142599999999
$ws.Range("L137").Value = 7413088.5
$ws.Range("M137").Value = -2933.142599999999
$ws.Range("N137").Value = -7418188.5

$ws.Range("H138").Value = 2651.7837
$ws.Range("I138").Value = 1778.5238
$ws.Range("J138").Value = 3797.9375
$ws.Range("K138").Value = 5335.5714
$ws.Range("L138").Value = 11393.8125
$ws.Range("M138").Value = -195.5713999999998
$ws.Range("N138").Value = -21673.8125

$ws.Range("H141").Value = 1054
$ws.Range("I141").Value = 1054
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3162
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2018
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1720.7778
$ws.Range("I22").Value = 1720.7778
$ws.Range("K22").Value = 1720.7778
$ws.Range("M22").Value = -1547.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2842.7144
$ws.Range("I12").Value = 2842.7144
$ws.Range("K12").Value = 2842.7144
$ws.Range("M12").Value = -2672.7144

$ws.Range("H31").Value = 6924.5
$ws.Range("I31").Value = 3587.4167
$ws.Range("K31").Value = 3587.4167
$ws.Range("M31").Value = -3292.4167

$ws.Range("H34").Value = 6924.5
$ws.Range("I34").Value = 3587.4167
$ws.Range("K34").Value = 3587.4167
$ws.Range("M34").Value = -3385.4167

$ws.Range("H42").Value = 25000
$ws.Range("J42").Value = 25000
$ws.Range("L42").Value = 25000
$ws.Range("N42").Value = -26186

$ws.Range("H53").Value = 57000
$ws.Range("J53").Value = 57000
$ws.Range("L53").Value = 57000
$ws.Range("N53").Value = -58214

$ws.Range("H107").Value = 2239.7273
$ws.Range("I107").Value = 2930
$ws.Range("K107").Value = 2930
$ws.Range("M107").Value = -1010

$ws.Range("H122").Value = 3551
$ws.Range("I122").Value = 2623
$ws.Range("K122").Value = 7869
$ws.Range("M122").Value = -5419

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I34").Value = 20000
$ws.Range("J34").Value = 65000
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 65000
$ws.Range("M34").Value = -19732
$ws.Range("N34").Value = -65536

$ws.Range("H70").Value = 5417.9165
$ws.Range("I70").Value = 5445.875
$ws.Range("J70").Value = 5362
$ws.Range("K70").Value = 5445.875
$ws.Range("L70").Value = 5362
$ws.Range("M70").Value = -5175.875
$ws.Range("N70").Value = -5902

$ws.Range("H73").Value = 5417.9165
$ws.Range("I73").Value = 5445.875
$ws.Range("J73").Value = 5362
$ws.Range("K73").Value = 5445.875
$ws.Range("L73").Value = 5362
$ws.Range("M73").Value = -4509.875
$ws.Range("N73").Value = -7234

$ws.Range("I76").Value = 20000
$ws.Range("J76").Value = 65000
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 65000
$ws.Range("M76").Value = -19685
$ws.Range("N76").Value = -65630

$ws.Range("I79").Value = 20000
$ws.Range("J79").Value = 65000
$ws.Range("K79").Value = 20000
$ws.Range("L79").Value = 65000
$ws.Range("M79").Value = -18908
$ws.Range("N79").Value = -67184

$ws.Range("H132").Value = 605982
$ws.Range("I132").Value = 1340950.1
$ws.Range("J132").Value = 4644.4546
$ws.Range("K132").Value = 4022850.3
$ws.Range("L132").Value = 13933.3638
$ws.Range("M132").Value = -4020320.3
$ws.Range("N132").Value = -18993.3638

$ws.Range("H136").Value = 23166.75
$ws.Range("J136").Value = 23166.75
$ws.Range("L136").Value = 69500.25
$ws.Range("N136").Value = -74600.25

$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 232.93333
$ws.Range("I55").Value = 200.42857
$ws.Range("J55").Value = 261.375
$ws.Range("K55").Value = 200.42857
$ws.Range("L55").Value = 261.375
$ws.Range("M55").Value = -27.42857000000001
$ws.Range("N55").Value = -607.375

$ws.Range("H82").Value = 1978.4
$ws.Range("I82").Value = 1332.8
$ws.Range("K82").Value = 1332.8
$ws.Range("M82").Value = -971.8

$ws.Range("H85").Value = 1978.4
$ws.Range("I85").Value = 1332.8
$ws.Range("K85").Value = 1332.8
$ws.Range("M85").Value = -84.79999999999995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 400000
$ws.Range("J15").Value = 400000
$ws.Range("L15").Value = 400000
$ws.Range("N15").Value = -400576

$ws.Range("H81").Value = 710.55554
$ws.Range("I81").Value = 650.7143
$ws.Range("J81").Value = 920
$ws.Range("K81").Value = 1301.4286
$ws.Range("L81").Value = 1840
$ws.Range("M81").Value = -240.4286
$ws.Range("N81").Value = -3962

$ws.Range("H84").Value = 710.55554
$ws.Range("I84").Value = 650.7143
$ws.Range("J84").Value = 920
$ws.Range("K84").Value = 6507.143
$ws.Range("L84").Value = 9200
$ws.Range("M84").Value = -1203.143
$ws.Range("N84").Value = -19808

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H141").Value = 96000
$ws.Range("J141").Value = 96000
$ws.Range("L141").Value = 96000
$ws.Range("N141").Value = -106360
